# Update odds values on the "Jogos da Semana" sheet to reflect the latest
# FlashScore data for 2024-10-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Jaguares de Cordoba vs Chico) updates
$ws.Range("G2").Value  = 1.67
$ws.Range("I2").Value  = 6
$ws.Range("J2").Value  = 2.38
$ws.Range("M2").Value  = 1.08
$ws.Range("N2").Value  = 8
$ws.Range("O2").Value  = 1.4
$ws.Range("P2").Value  = 2.75
$ws.Range("X2").Value  = 6.5
$ws.Range("Z2").Value  = 12
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 29
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 67
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 29
$ws.Range("AR2").Value = 51
$ws.Range("AW2").Value = 7
$ws.Range("AX2").Value = 34

# Row 3 (Liverpool M. vs Wanderers) updates
$ws.Range("M3").Value = 1.03
$ws.Range("O3").Value = 1.27
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.85
